# Set default warm offset to 550
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biases")

# Rows 34-37 correspond to DET0_OFFSET, DET1_OFFSET, DET2_OFFSET, DET3_OFFSET
# Columns B (2) through BD (56) hold the per-module values, currently 0.
$ws.Range("B34:BD37").Value = 550
